$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.987.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.364.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.363.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.468'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  -4.10%  '
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.938.71'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('E16').Value = '  -5.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.364.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.107.59'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.553'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.498.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.21%  '
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('E31').Value = '  +3.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -6.44%  '
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.82'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0751'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.78%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.768'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  -2.32%  '
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.349.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('E51').Value = '  +1.83%  '
